# Pulse_BOM.xlsx update:
#  - Rename the sheet tab from "Pulse" to "BOM"
#  - Fix the bi-color LED naming: "LED 3mm Flat Bicolor" -> "LED 3mm Dome Bicolor"
#  - Leave the active cell/selection on C18 (where the renamed part lives)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "BOM"

# Correct the LED description in the BOM table (row 18, Package column).
$ws.Range("C18").Value = "LED 3mm Dome Bicolor"

# Match the saved selection/active cell.
$ws.Range("C18").Select()
